$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2573.7144
$ws.Range("J58").Value = 12017
$ws.Range("L58").Value = 36051
$ws.Range("N58").Value = -36351
$ws.Range("H64").Value = 4886.4
$ws.Range("I64").Value = 4251.1
$ws.Range("K64").Value = 4251.1
$ws.Range("M64").Value = -4003.1
$ws.Range("H67").Value = 4886.4
$ws.Range("I67").Value = 4251.1
$ws.Range("K67").Value = 4251.1
$ws.Range("M67").Value = -3393.1
$ws.Range("H69").Value = 16723.334
$ws.Range("J69").Value = 4984
$ws.Range("L69").Value = 14952
$ws.Range("N69").Value = -16700
$ws.Range("H72").Value = 16723.334
$ws.Range("J72").Value = 4984
$ws.Range("L72").Value = 44856
$ws.Range("N72").Value = -53592
$ws.Range("H74").Value = 5753.5386
$ws.Range("I74").Value = 3632.8333
$ws.Range("K74").Value = 3632.8333
$ws.Range("M74").Value = -2696.8333
$ws.Range("H76").Value = 4999.8335
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H77").Value = 5753.5386
$ws.Range("I77").Value = 3632.8333
$ws.Range("K77").Value = 18164.1665
$ws.Range("M77").Value = -13484.1665
$ws.Range("H79").Value = 4999.8335
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H107").Value = 1016.125
$ws.Range("I107").Value = 1008.8947
$ws.Range("J107").Value = 1043.6
$ws.Range("K107").Value = 1008.8947
$ws.Range("L107").Value = 1043.6
$ws.Range("M107").Value = 911.1053000000001
$ws.Range("N107").Value = -4883.6
$ws.Range("H118").Value = 576.25
$ws.Range("I118").Value = 576.25
$ws.Range("K118").Value = 1728.75
$ws.Range("M118").Value = -71.75
$ws.Range("H138").Value = 3767
$ws.Range("J138").Value = 5993
$ws.Range("L138").Value = 17979
$ws.Range("N138").Value = -28259
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9994.4375
$ws.Range("J2").Value = 4672.6665
$ws.Range("L2").Value = 4672.6665
$ws.Range("N2").Value = -4898.6665
$ws.Range("H45").Value = 1720.25
$ws.Range("I45").Value = 1649
$ws.Range("K45").Value = 1649
$ws.Range("M45").Value = -1272
$ws.Range("H74").Value = 3562889.5
$ws.Range("I74").Value = 4630497
$ws.Range("J74").Value = 3088397
$ws.Range("K74").Value = 4630497
$ws.Range("L74").Value = 3088397
$ws.Range("M74").Value = -4629623
$ws.Range("N74").Value = -3090145
$ws.Range("H77").Value = 3562889.5
$ws.Range("I77").Value = 4630497
$ws.Range("J77").Value = 3088397
$ws.Range("K77").Value = 23152485
$ws.Range("L77").Value = 15441985
$ws.Range("M77").Value = -23148117
$ws.Range("N77").Value = -15450721
$ws.Range("H102").Value = 1483.375
$ws.Range("I102").Value = 1429
$ws.Range("J102").Value = 2299
$ws.Range("K102").Value = 1429
$ws.Range("L102").Value = 2299
$ws.Range("M102").Value = 193
$ws.Range("N102").Value = -5543
$ws.Range("H110").Value = 7230
$ws.Range("I110").Value = 6757.1875
$ws.Range("K110").Value = 6757.1875
$ws.Range("M110").Value = -4712.1875
$ws.Range("H116").Value = 9994.4375
$ws.Range("J116").Value = 4672.6665
$ws.Range("L116").Value = 4672.6665
$ws.Range("N116").Value = -9260.666499999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9994.4375
$ws.Range("J3").Value = 4672.6665
$ws.Range("L3").Value = 4672.6665
$ws.Range("N3").Value = -4900.6665
$ws.Range("H20").Value = 2145.4285
$ws.Range("I20").Value = 1687.2
$ws.Range("J20").Value = 3291
$ws.Range("K20").Value = 1687.2
$ws.Range("L20").Value = 3291
$ws.Range("M20").Value = -1440.2
$ws.Range("N20").Value = -3785
$ws.Range("H86").Value = 1591.5
$ws.Range("I86").Value = 1283.3334
$ws.Range("J86").Value = 1899.6666
$ws.Range("K86").Value = 1283.3334
$ws.Range("L86").Value = 1899.6666
$ws.Range("M86").Value = -160.3334
$ws.Range("N86").Value = -4145.6666
$ws.Range("H89").Value = 1591.5
$ws.Range("I89").Value = 1283.3334
$ws.Range("J89").Value = 1899.6666
$ws.Range("K89").Value = 6416.666999999999
$ws.Range("L89").Value = 9498.333000000001
$ws.Range("M89").Value = -800.6669999999995
$ws.Range("N89").Value = -20730.333
$ws.Range("H99").Value = 2518.611
$ws.Range("I99").Value = 2412.8
$ws.Range("K99").Value = 2412.8
$ws.Range("M99").Value = -914.8000000000002
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3208.9033
$ws.Range("I22").Value = 4507.154
$ws.Range("K22").Value = 4507.154
$ws.Range("M22").Value = -4157.154
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 3424.2856
$ws.Range("J69").Value = 3424.2856
$ws.Range("L69").Value = 10272.8568
$ws.Range("N69").Value = -11894.8568
$ws.Range("H72").Value = 3424.2856
$ws.Range("J72").Value = 3424.2856
$ws.Range("L72").Value = 30818.5704
$ws.Range("N72").Value = -38930.5704
$ws.Range("H80").Value = 4500
$ws.Range("J80").Value = 4500
$ws.Range("L80").Value = 13500
$ws.Range("N80").Value = -15372
$ws.Range("H83").Value = 4500
$ws.Range("J83").Value = 4500
$ws.Range("L83").Value = 40500
$ws.Range("N83").Value = -49860
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8316.166999999999
$ws.Range("I70").Value = 8736
$ws.Range("K70").Value = 8736
$ws.Range("M70").Value = -8466
$ws.Range("H73").Value = 8316.166999999999
$ws.Range("I73").Value = 8736
$ws.Range("K73").Value = 8736
$ws.Range("M73").Value = -7800
$ws.Range("H80").Value = 8904.637000000001
$ws.Range("I80").Value = 2517.6667
$ws.Range("K80").Value = 2517.6667
$ws.Range("M80").Value = -1519.6667
$ws.Range("H83").Value = 8904.637000000001
$ws.Range("I83").Value = 2517.6667
$ws.Range("K83").Value = 12588.3335
$ws.Range("M83").Value = -7596.333500000001
$ws.Range("H97").Value = 1016.9091
$ws.Range("I97").Value = 848.44446
$ws.Range("K97").Value = 848.44446
$ws.Range("M97").Value = -352.44446
$ws.Range("H102").Value = 3480.0833
$ws.Range("I102").Value = 3170.5
$ws.Range("K102").Value = 3170.5
$ws.Range("M102").Value = -1548.5
$ws.Range("H107").Value = 670.7059
$ws.Range("I107").Value = 430.58334
$ws.Range("J107").Value = 1247
$ws.Range("K107").Value = 430.58334
$ws.Range("L107").Value = 1247
$ws.Range("M107").Value = 1489.41666
$ws.Range("N107").Value = -5087
$ws.Range("H132").Value = 2969.1428
$ws.Range("I132").Value = 2958.8
$ws.Range("J132").Value = 2995
$ws.Range("K132").Value = 8876.400000000001
$ws.Range("L132").Value = 8985
$ws.Range("M132").Value = -6346.400000000001
$ws.Range("N132").Value = -14045
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9791.154
$ws.Range("I7").Value = 9899.416999999999
$ws.Range("J7").Value = 8492
$ws.Range("K7").Value = 9899.416999999999
$ws.Range("L7").Value = 8492
$ws.Range("M7").Value = -9787.416999999999
$ws.Range("N7").Value = -8716
$ws.Range("H16").Value = 1412.9412
$ws.Range("I16").Value = 1412.9412
$ws.Range("K16").Value = 1412.9412
$ws.Range("M16").Value = -1242.9412
$ws.Range("H19").Value = 821.4
$ws.Range("J19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("N19").Value = -3340
$ws.Range("H126").Value = 9791.154
$ws.Range("I126").Value = 9899.416999999999
$ws.Range("J126").Value = 8492
$ws.Range("K126").Value = 29698.251
$ws.Range("L126").Value = 25476
$ws.Range("M126").Value = -27228.251
$ws.Range("N126").Value = -30416
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 160692.2
$ws.Range("I74").Value = 6998.5
$ws.Range("J74").Value = 263154.66
$ws.Range("K74").Value = 6998.5
$ws.Range("L74").Value = 263154.66
$ws.Range("M74").Value = -6062.5
$ws.Range("N74").Value = -265026.66
$ws.Range("H77").Value = 160692.2
$ws.Range("I77").Value = 6998.5
$ws.Range("J77").Value = 263154.66
$ws.Range("K77").Value = 20995.5
$ws.Range("L77").Value = 789463.98
$ws.Range("M77").Value = -16315.5
$ws.Range("N77").Value = -798823.98
$ws.Range("H107").Value = 1957.4375
$ws.Range("I107").Value = 688.5
$ws.Range("J107").Value = 4072.3333
$ws.Range("K107").Value = 2065.5
$ws.Range("L107").Value = 12216.9999
$ws.Range("M107").Value = -145.5
$ws.Range("N107").Value = -16056.9999
